$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.495
$ws.Range("A3").Value = -21.611
$ws.Range("E4").Value = 13.073
$ws.Range("B5").Value = 6.339
$ws.Range("C5").Value = -11.998
$ws.Range("D7").Value = -7.105000000000001
$ws.Range("E7").Value = 13.398
$ws.Range("C9").Value = -12.394
$ws.Range("C11").Value = -12.089
$ws.Range("D11").Value = -7.74
$ws.Range("A14").Value = -20.945
$ws.Range("A16").Value = -21.363
$ws.Range("B16").Value = 6.007999999999999
$ws.Range("C17").Value = -11.61
$ws.Range("D19").Value = -7.838000000000001
$ws.Range("A21").Value = -22.003
$ws.Range("C21").Value = -12.18
$ws.Range("D21").Value = -7.941
$ws.Range("A23").Value = -21.218
$ws.Range("E23").Value = 13.741
$ws.Range("A25").Value = -22.195
